$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
$ps.Draft = $false
$ps.BlackAndWhite = $false
$ps.PrintQuality = 600
